$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite full table A1:K37 to match target layout (rows/labels shifted; Fossil Gases & Fossil Liquids rows inserted per year)

# Row 1
$ws.Range("A1").Value2 = "FuelGroup"
$ws.Range("B1").Value2 = "Year"
$ws.Range("C1").Value2 = "Iron & steel"
$ws.Range("D1").Value2 = "Chemicals"
$ws.Range("E1").Value2 = "Non-metallic minerals"
$ws.Range("F1").Value2 = "Pass Road"
$ws.Range("G1").Value2 = "Pass Rail"
$ws.Range("H1").Value2 = "Pass Aviation"
$ws.Range("I1").Value2 = "Freight Road"
$ws.Range("J1").Value2 = "Freight Rail"
$ws.Range("K1").Value2 = "Maritime"

# Row 2
$ws.Range("A2").Value2 = "Hydrogen"
$ws.Range("B2").Value2 = 2030
$ws.Range("C2").Value2 = $null
$ws.Range("D2").Value2 = $null
$ws.Range("E2").Value2 = $null
$ws.Range("F2").Value2 = 0.0036479357571302
$ws.Range("G2").Value2 = $null
$ws.Range("H2").Value2 = 0.00000001335177568378488
$ws.Range("I2").Value2 = 0.0016425518663523
$ws.Range("J2").Value2 = $null
$ws.Range("K2").Value2 = $null

# Row 3
$ws.Range("A3").Value2 = "Methanol"
$ws.Range("B3").Value2 = 2030
$ws.Range("C3").Value2 = $null
$ws.Range("D3").Value2 = 0.0001667221954613387
$ws.Range("E3").Value2 = $null
$ws.Range("F3").Value2 = $null
$ws.Range("G3").Value2 = $null
$ws.Range("H3").Value2 = $null
$ws.Range("I3").Value2 = $null
$ws.Range("J3").Value2 = $null
$ws.Range("K3").Value2 = $null

# Row 4
$ws.Range("A4").Value2 = "Ammonia"
$ws.Range("B4").Value2 = 2030
$ws.Range("C4").Value2 = $null
$ws.Range("D4").Value2 = 0.01272982868714951
$ws.Range("E4").Value2 = $null
$ws.Range("F4").Value2 = $null
$ws.Range("G4").Value2 = $null
$ws.Range("H4").Value2 = $null
$ws.Range("I4").Value2 = $null
$ws.Range("J4").Value2 = $null
$ws.Range("K4").Value2 = $null

# Row 5
$ws.Range("A5").Value2 = "Synthetic Gases"
$ws.Range("B5").Value2 = 2030
$ws.Range("C5").Value2 = $null
$ws.Range("D5").Value2 = $null
$ws.Range("E5").Value2 = $null
$ws.Range("F5").Value2 = $null
$ws.Range("G5").Value2 = $null
$ws.Range("H5").Value2 = $null
$ws.Range("I5").Value2 = $null
$ws.Range("J5").Value2 = $null
$ws.Range("K5").Value2 = $null

# Row 6
$ws.Range("A6").Value2 = "Biogenic Gases"
$ws.Range("B6").Value2 = 2030
$ws.Range("C6").Value2 = $null
$ws.Range("D6").Value2 = $null
$ws.Range("E6").Value2 = 0.001517032094599047
$ws.Range("F6").Value2 = 0.0011196654984368
$ws.Range("G6").Value2 = $null
$ws.Range("H6").Value2 = $null
$ws.Range("I6").Value2 = 0.0002953533760394
$ws.Range("J6").Value2 = $null
$ws.Range("K6").Value2 = $null

# Row 7
$ws.Range("A7").Value2 = "Fossil Gases"
$ws.Range("B7").Value2 = 2030
$ws.Range("C7").Value2 = $null
$ws.Range("D7").Value2 = $null
$ws.Range("E7").Value2 = $null
$ws.Range("F7").Value2 = 0.0140462882103178
$ws.Range("G7").Value2 = $null
$ws.Range("H7").Value2 = $null
$ws.Range("I7").Value2 = 0.0012192777718049
$ws.Range("J7").Value2 = $null
$ws.Range("K7").Value2 = $null

# Row 8
$ws.Range("A8").Value2 = "Synthetic Liquids"
$ws.Range("B8").Value2 = 2030
$ws.Range("C8").Value2 = $null
$ws.Range("D8").Value2 = $null
$ws.Range("E8").Value2 = $null
$ws.Range("F8").Value2 = $null
$ws.Range("G8").Value2 = $null
$ws.Range("H8").Value2 = $null
$ws.Range("I8").Value2 = $null
$ws.Range("J8").Value2 = $null
$ws.Range("K8").Value2 = $null

# Row 9
$ws.Range("A9").Value2 = "Biogenic Liquids"
$ws.Range("B9").Value2 = 2030
$ws.Range("C9").Value2 = $null
$ws.Range("D9").Value2 = $null
$ws.Range("E9").Value2 = $null
$ws.Range("F9").Value2 = 0.0696096341435164
$ws.Range("G9").Value2 = 0.0005003989726162
$ws.Range("H9").Value2 = 0.0147095845978992
$ws.Range("I9").Value2 = 0.0576730690117931
$ws.Range("J9").Value2 = 0.0004268491355289
$ws.Range("K9").Value2 = 0.0250105603733093

# Row 10
$ws.Range("A10").Value2 = "Fossil Liquids"
$ws.Range("B10").Value2 = 2030
$ws.Range("C10").Value2 = $null
$ws.Range("D10").Value2 = $null
$ws.Range("E10").Value2 = $null
$ws.Range("F10").Value2 = 0.6962485521627411
$ws.Range("G10").Value2 = 0.0035241632245806
$ws.Range("H10").Value2 = 0.1305643970668376
$ws.Range("I10").Value2 = 0.3639310136313932
$ws.Range("J10").Value2 = 0.0025837187153698
$ws.Range("K10").Value2 = 0.2348176943229761

# Row 11
$ws.Range("A11").Value2 = "Biomass [Solid]"
$ws.Range("B11").Value2 = 2030
$ws.Range("C11").Value2 = $null
$ws.Range("D11").Value2 = $null
$ws.Range("E11").Value2 = 0.002537305967593918
$ws.Range("F11").Value2 = $null
$ws.Range("G11").Value2 = $null
$ws.Range("H11").Value2 = $null
$ws.Range("I11").Value2 = $null
$ws.Range("J11").Value2 = $null
$ws.Range("K11").Value2 = $null

# Row 12
$ws.Range("A12").Value2 = "Renewable Energy Carrier"
$ws.Range("B12").Value2 = 2030
$ws.Range("C12").Value2 = $null
$ws.Range("D12").Value2 = $null
$ws.Range("E12").Value2 = 0.0006271917820834629
$ws.Range("F12").Value2 = $null
$ws.Range("G12").Value2 = $null
$ws.Range("H12").Value2 = $null
$ws.Range("I12").Value2 = $null
$ws.Range("J12").Value2 = $null
$ws.Range("K12").Value2 = $null

# Row 13
$ws.Range("A13").Value2 = "Overall Demand"
$ws.Range("B13").Value2 = 2030
$ws.Range("C13").Value2 = $null
$ws.Range("D13").Value2 = 0.01289655088261085
$ws.Range("E13").Value2 = 0.004681529844276429
$ws.Range("F13").Value2 = 0.7846720757721423
$ws.Range("G13").Value2 = 0.0040245621971968
$ws.Range("H13").Value2 = 0.1452739950165125
$ws.Range("I13").Value2 = 0.424761265657383
$ws.Range("J13").Value2 = 0.0030105678508987
$ws.Range("K13").Value2 = 0.2598282546962854

# Row 14
$ws.Range("A14").Value2 = "Hydrogen"
$ws.Range("B14").Value2 = 2040
$ws.Range("C14").Value2 = $null
$ws.Range("D14").Value2 = $null
$ws.Range("E14").Value2 = $null
$ws.Range("F14").Value2 = 0.0175253502218624
$ws.Range("G14").Value2 = $null
$ws.Range("H14").Value2 = 0.000001117691096132891
$ws.Range("I14").Value2 = 0.0023865969783612
$ws.Range("J14").Value2 = $null
$ws.Range("K14").Value2 = $null

# Row 15
$ws.Range("A15").Value2 = "Methanol"
$ws.Range("B15").Value2 = 2040
$ws.Range("C15").Value2 = $null
$ws.Range("D15").Value2 = 0.0001819930503096124
$ws.Range("E15").Value2 = $null
$ws.Range("F15").Value2 = $null
$ws.Range("G15").Value2 = $null
$ws.Range("H15").Value2 = $null
$ws.Range("I15").Value2 = $null
$ws.Range("J15").Value2 = $null
$ws.Range("K15").Value2 = $null

# Row 16
$ws.Range("A16").Value2 = "Ammonia"
$ws.Range("B16").Value2 = 2040
$ws.Range("C16").Value2 = $null
$ws.Range("D16").Value2 = 0.01389581240987423
$ws.Range("E16").Value2 = $null
$ws.Range("F16").Value2 = $null
$ws.Range("G16").Value2 = $null
$ws.Range("H16").Value2 = $null
$ws.Range("I16").Value2 = $null
$ws.Range("J16").Value2 = $null
$ws.Range("K16").Value2 = $null

# Row 17
$ws.Range("A17").Value2 = "Synthetic Gases"
$ws.Range("B17").Value2 = 2040
$ws.Range("C17").Value2 = $null
$ws.Range("D17").Value2 = $null
$ws.Range("E17").Value2 = $null
$ws.Range("F17").Value2 = 0.000000006421190540173386
$ws.Range("G17").Value2 = $null
$ws.Range("H17").Value2 = $null
$ws.Range("I17").Value2 = 0.0000000007404392471644982
$ws.Range("J17").Value2 = $null
$ws.Range("K17").Value2 = $null

# Row 18
$ws.Range("A18").Value2 = "Biogenic Gases"
$ws.Range("B18").Value2 = 2040
$ws.Range("C18").Value2 = $null
$ws.Range("D18").Value2 = $null
$ws.Range("E18").Value2 = 0.002055309681044841
$ws.Range("F18").Value2 = 0.0013764905176361
$ws.Range("G18").Value2 = $null
$ws.Range("H18").Value2 = $null
$ws.Range("I18").Value2 = 0.000509620242557
$ws.Range("J18").Value2 = $null
$ws.Range("K18").Value2 = $null

# Row 19
$ws.Range("A19").Value2 = "Fossil Gases"
$ws.Range("B19").Value2 = 2040
$ws.Range("C19").Value2 = $null
$ws.Range("D19").Value2 = $null
$ws.Range("E19").Value2 = $null
$ws.Range("F19").Value2 = 0.007613812729963001
$ws.Range("G19").Value2 = $null
$ws.Range("H19").Value2 = $null
$ws.Range("I19").Value2 = 0.0012931202386089
$ws.Range("J19").Value2 = $null
$ws.Range("K19").Value2 = $null

# Row 20
$ws.Range("A20").Value2 = "Synthetic Liquids"
$ws.Range("B20").Value2 = 2040
$ws.Range("C20").Value2 = $null
$ws.Range("D20").Value2 = $null
$ws.Range("E20").Value2 = $null
$ws.Range("F20").Value2 = $null
$ws.Range("G20").Value2 = $null
$ws.Range("H20").Value2 = $null
$ws.Range("I20").Value2 = $null
$ws.Range("J20").Value2 = $null
$ws.Range("K20").Value2 = $null

# Row 21
$ws.Range("A21").Value2 = "Biogenic Liquids"
$ws.Range("B21").Value2 = 2040
$ws.Range("C21").Value2 = $null
$ws.Range("D21").Value2 = $null
$ws.Range("E21").Value2 = $null
$ws.Range("F21").Value2 = 0.0325488940884981
$ws.Range("G21").Value2 = 0.0008157328785995
$ws.Range("H21").Value2 = 0.0191863575218352
$ws.Range("I21").Value2 = 0.0380820558157449
$ws.Range("J21").Value2 = 0.0005185524227517
$ws.Range("K21").Value2 = 0.0289957981166524

# Row 22
$ws.Range("A22").Value2 = "Fossil Liquids"
$ws.Range("B22").Value2 = 2040
$ws.Range("C22").Value2 = $null
$ws.Range("D22").Value2 = $null
$ws.Range("E22").Value2 = $null
$ws.Range("F22").Value2 = 0.2148972398498224
$ws.Range("G22").Value2 = 0.0037886713998633
$ws.Range("H22").Value2 = 0.1236661338043652
$ws.Range("I22").Value2 = 0.1626073741385843
$ws.Range("J22").Value2 = 0.0022928192073582
$ws.Range("K22").Value2 = 0.2253953088211194

# Row 23
$ws.Range("A23").Value2 = "Biomass [Solid]"
$ws.Range("B23").Value2 = 2040
$ws.Range("C23").Value2 = $null
$ws.Range("D23").Value2 = $null
$ws.Range("E23").Value2 = $null
$ws.Range("F23").Value2 = $null
$ws.Range("G23").Value2 = $null
$ws.Range("H23").Value2 = $null
$ws.Range("I23").Value2 = $null
$ws.Range("J23").Value2 = $null
$ws.Range("K23").Value2 = $null

# Row 24
$ws.Range("A24").Value2 = "Renewable Energy Carrier"
$ws.Range("B24").Value2 = 2040
$ws.Range("C24").Value2 = $null
$ws.Range("D24").Value2 = $null
$ws.Range("E24").Value2 = $null
$ws.Range("F24").Value2 = $null
$ws.Range("G24").Value2 = $null
$ws.Range("H24").Value2 = $null
$ws.Range("I24").Value2 = $null
$ws.Range("J24").Value2 = $null
$ws.Range("K24").Value2 = $null

# Row 25
$ws.Range("A25").Value2 = "Overall Demand"
$ws.Range("B25").Value2 = 2040
$ws.Range("C25").Value2 = $null
$ws.Range("D25").Value2 = 0.01407780546018384
$ws.Range("E25").Value2 = 0.002055309681044841
$ws.Range("F25").Value2 = 0.2739617938289725
$ws.Range("G25").Value2 = 0.0046044042784628
$ws.Range("H25").Value2 = 0.1428536090172965
$ws.Range("I25").Value2 = 0.2048787681542955
$ws.Range("J25").Value2 = 0.0028113716301099
$ws.Range("K25").Value2 = 0.2543911069377718

# Row 26
$ws.Range("A26").Value2 = "Hydrogen"
$ws.Range("B26").Value2 = 2050
$ws.Range("C26").Value2 = $null
$ws.Range("D26").Value2 = $null
$ws.Range("E26").Value2 = $null
$ws.Range("F26").Value2 = 0.0243107274205465
$ws.Range("G26").Value2 = $null
$ws.Range("H26").Value2 = 0.000001894396379722037
$ws.Range("I26").Value2 = 0.003804105309536
$ws.Range("J26").Value2 = $null
$ws.Range("K26").Value2 = $null

# Row 27
$ws.Range("A27").Value2 = "Methanol"
$ws.Range("B27").Value2 = 2050
$ws.Range("C27").Value2 = $null
$ws.Range("D27").Value2 = 0.000195852152862953
$ws.Range("E27").Value2 = $null
$ws.Range("F27").Value2 = $null
$ws.Range("G27").Value2 = $null
$ws.Range("H27").Value2 = $null
$ws.Range("I27").Value2 = $null
$ws.Range("J27").Value2 = $null
$ws.Range("K27").Value2 = $null

# Row 28
$ws.Range("A28").Value2 = "Ammonia"
$ws.Range("B28").Value2 = 2050
$ws.Range("C28").Value2 = $null
$ws.Range("D28").Value2 = 0.01495400385687069
$ws.Range("E28").Value2 = $null
$ws.Range("F28").Value2 = $null
$ws.Range("G28").Value2 = $null
$ws.Range("H28").Value2 = $null
$ws.Range("I28").Value2 = $null
$ws.Range("J28").Value2 = $null
$ws.Range("K28").Value2 = $null

# Row 29
$ws.Range("A29").Value2 = "Synthetic Gases"
$ws.Range("B29").Value2 = 2050
$ws.Range("C29").Value2 = $null
$ws.Range("D29").Value2 = $null
$ws.Range("E29").Value2 = $null
$ws.Range("F29").Value2 = 0.00000006040957825554191
$ws.Range("G29").Value2 = $null
$ws.Range("H29").Value2 = $null
$ws.Range("I29").Value2 = 0.00000001862117208524179
$ws.Range("J29").Value2 = $null
$ws.Range("K29").Value2 = $null

# Row 30
$ws.Range("A30").Value2 = "Biogenic Gases"
$ws.Range("B30").Value2 = 2050
$ws.Range("C30").Value2 = $null
$ws.Range("D30").Value2 = $null
$ws.Range("E30").Value2 = 0.003784923913590072
$ws.Range("F30").Value2 = 0.0002362118286557918
$ws.Range("G30").Value2 = $null
$ws.Range("H30").Value2 = $null
$ws.Range("I30").Value2 = 0.0001464881637984659
$ws.Range("J30").Value2 = $null
$ws.Range("K30").Value2 = $null

# Row 31
$ws.Range("A31").Value2 = "Fossil Gases"
$ws.Range("B31").Value2 = 2050
$ws.Range("C31").Value2 = $null
$ws.Range("D31").Value2 = $null
$ws.Range("E31").Value2 = $null
$ws.Range("F31").Value2 = 0.0004771160107088
$ws.Range("G31").Value2 = $null
$ws.Range("H31").Value2 = $null
$ws.Range("I31").Value2 = 0.0004831337696037001
$ws.Range("J31").Value2 = $null
$ws.Range("K31").Value2 = $null

# Row 32
$ws.Range("A32").Value2 = "Synthetic Liquids"
$ws.Range("B32").Value2 = 2050
$ws.Range("C32").Value2 = $null
$ws.Range("D32").Value2 = $null
$ws.Range("E32").Value2 = $null
$ws.Range("F32").Value2 = 0.0000000006626646519652098
$ws.Range("G32").Value2 = 0.00000000003003903291544851
$ws.Range("H32").Value2 = 0.000000000634537903847718
$ws.Range("I32").Value2 = 0.0000000006339460051053089
$ws.Range("J32").Value2 = 0.000000000003551475162335466
$ws.Range("K32").Value2 = 0.000000002488598384219159

# Row 33
$ws.Range("A33").Value2 = "Biogenic Liquids"
$ws.Range("B33").Value2 = 2050
$ws.Range("C33").Value2 = $null
$ws.Range("D33").Value2 = $null
$ws.Range("E33").Value2 = $null
$ws.Range("F33").Value2 = 0.0067572942771184
$ws.Range("G33").Value2 = 0.0014575714346724
$ws.Range("H33").Value2 = 0.0267308908885776
$ws.Range("I33").Value2 = 0.0098166842558465
$ws.Range("J33").Value2 = 0.0006669405310802
$ws.Range("K33").Value2 = 0.0409729797312409

# Row 34
$ws.Range("A34").Value2 = "Fossil Liquids"
$ws.Range("B34").Value2 = 2050
$ws.Range("C34").Value2 = $null
$ws.Range("D34").Value2 = $null
$ws.Range("E34").Value2 = $null
$ws.Range("F34").Value2 = 0.0388318373000604
$ws.Range("G34").Value2 = 0.0034172211511031
$ws.Range("H34").Value2 = 0.106560884069033
$ws.Range("I34").Value2 = 0.0289959298007036
$ws.Range("J34").Value2 = 0.0019735461906695
$ws.Range("K34").Value2 = 0.2074284198054586

# Row 35
$ws.Range("A35").Value2 = "Biomass [Solid]"
$ws.Range("B35").Value2 = 2050
$ws.Range("C35").Value2 = $null
$ws.Range("D35").Value2 = $null
$ws.Range("E35").Value2 = $null
$ws.Range("F35").Value2 = $null
$ws.Range("G35").Value2 = $null
$ws.Range("H35").Value2 = $null
$ws.Range("I35").Value2 = $null
$ws.Range("J35").Value2 = $null
$ws.Range("K35").Value2 = $null

# Row 36
$ws.Range("A36").Value2 = "Renewable Energy Carrier"
$ws.Range("B36").Value2 = 2050
$ws.Range("C36").Value2 = $null
$ws.Range("D36").Value2 = $null
$ws.Range("E36").Value2 = $null
$ws.Range("F36").Value2 = $null
$ws.Range("G36").Value2 = $null
$ws.Range("H36").Value2 = $null
$ws.Range("I36").Value2 = $null
$ws.Range("J36").Value2 = $null
$ws.Range("K36").Value2 = $null

# Row 37
$ws.Range("A37").Value2 = "Overall Demand"
$ws.Range("B37").Value2 = 2050
$ws.Range("C37").Value2 = $null
$ws.Range("D37").Value2 = 0.01514985600973365
$ws.Range("E37").Value2 = 0.003784923913590072
$ws.Range("F37").Value2 = 0.07061324790933281
$ws.Range("G37").Value2 = 0.004874792615814533
$ws.Range("H37").Value2 = 0.1332936699885282
$ws.Range("I37").Value2 = 0.04324636055460636
$ws.Range("J37").Value2 = 0.002640486725301175
$ws.Range("K37").Value2 = 0.2484014020252979
